$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (CO1 / CustomNBC / 30%) - updated metrics
$ws.Range("D8").Value = 0.39
$ws.Range("E8").Value = 0.69
$ws.Range("F8").Value = 0.47
$ws.Range("G8").Value = 0.56
$ws.Range("H8").Value = 0.63

# Row 9 (CO1 / 50%) - renamed from Mothur to CustomNBC, updated Precision
$ws.Range("B9").Value = "CustomNBC"
$ws.Range("E9").Value = 0.67
